$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 220
$ws.Range("I6").Value = 237.5
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 712.5
$ws.Range("L6").Value = 450
$ws.Range("M6").Value = -600.5
$ws.Range("N6").Value = -674
$ws.Range("H76").Value = 2688.6
$ws.Range("I76").Value = 2610.75
$ws.Range("K76").Value = 2610.75
$ws.Range("M76").Value = -2295.75
$ws.Range("H79").Value = 2688.6
$ws.Range("I79").Value = 2610.75
$ws.Range("K79").Value = 2610.75
$ws.Range("M79").Value = -1518.75
$ws.Range("H82").Value = 333
$ws.Range("I82").Value = 333
$ws.Range("K82").Value = 999
$ws.Range("M82").Value = -593
$ws.Range("H85").Value = 333
$ws.Range("I85").Value = 333
$ws.Range("K85").Value = 999
$ws.Range("M85").Value = 405
$ws.Range("H98").Value = 1295
$ws.Range("I98").Value = 1318.64
$ws.Range("J98").Value = 999.5
$ws.Range("K98").Value = 1318.64
$ws.Range("L98").Value = 999.5
$ws.Range("M98").Value = 179.3599999999999
$ws.Range("N98").Value = -3995.5
$ws.Range("H101").Value = 137
$ws.Range("I101").Value = 137
$ws.Range("K101").Value = 411
$ws.Range("M101").Value = 1211
$ws.Range("H122").Value = 1295
$ws.Range("I122").Value = 1318.64
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 3955.92
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -1505.92
$ws.Range("N122").Value = -7898.5
$ws.Range("H129").Value = 601.1429000000001
$ws.Range("I129").Value = 601.1429000000001
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1803.4287
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3196.5713
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 277.375
$ws.Range("I4").Value = 169.83333
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 169.83333
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -53.83332999999999
$ws.Range("N4").Value = -832
$ws.Range("H32").Value = 3012.3767
$ws.Range("I32").Value = 2238.4062
$ws.Range("J32").Value = 12919.2
$ws.Range("K32").Value = 2238.4062
$ws.Range("L32").Value = 12919.2
$ws.Range("M32").Value = -1951.4062
$ws.Range("N32").Value = -13493.2
$ws.Range("H61").Value = 2129
$ws.Range("I61").Value = 2129
$ws.Range("K61").Value = 2129
$ws.Range("M61").Value = -1917
$ws.Range("H63").Value = 2600
$ws.Range("I63").Value = 2600
$ws.Range("K63").Value = 2600
$ws.Range("M63").Value = -1914
$ws.Range("H66").Value = 2600
$ws.Range("I66").Value = 2600
$ws.Range("K66").Value = 13000
$ws.Range("M66").Value = -9568
$ws.Range("H88").Value = 1359.625
$ws.Range("I88").Value = 774.5
$ws.Range("K88").Value = 774.5
$ws.Range("M88").Value = -368.5
$ws.Range("H91").Value = 1359.625
$ws.Range("I91").Value = 774.5
$ws.Range("K91").Value = 774.5
$ws.Range("M91").Value = 629.5
$ws.Range("H136").Value = 2129
$ws.Range("I136").Value = 2129
$ws.Range("K136").Value = 6387
$ws.Range("M136").Value = -3837
$ws.Range("H141").Value = 48331.668
$ws.Range("J141").Value = 48331.668
$ws.Range("L141").Value = 48331.668
$ws.Range("N141").Value = -58691.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 7000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H134").Value = 1458.8572
$ws.Range("I134").Value = 1281.8
$ws.Range("K134").Value = 3845.4
$ws.Range("M134").Value = -1310.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10333.967
$ws.Range("I99").Value = 6039.0586
$ws.Range("K99").Value = 6039.0586
$ws.Range("M99").Value = -4541.0586
$ws.Range("H126").Value = 10333.967
$ws.Range("I126").Value = 6039.0586
$ws.Range("K126").Value = 18117.1758
$ws.Range("M126").Value = -15647.1758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 14285854
$ws.Range("I7").Value = 16666704
$ws.Range("K7").Value = 50000112
$ws.Range("M7").Value = -50000000
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 40
$ws.Range("K16").Value = 120
$ws.Range("M16").Value = 53
$ws.Range("H19").Value = 1900
$ws.Range("I19").Value = 1900
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5700
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5526
$ws.Range("N19").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 2999
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10869
$ws.Range("H83").Value = 2999
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36351
$ws.Range("H92").Value = 243.63637
$ws.Range("J92").Value = 218.2
$ws.Range("L92").Value = 654.5999999999999
$ws.Range("N92").Value = -3150.6
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 15000000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H43").Value = 318400.5
$ws.Range("I43").Value = 3727.5
$ws.Range("K43").Value = 3727.5
$ws.Range("M43").Value = -3534.5
$ws.Range("H53").Value = 12682
$ws.Range("I53").Value = 12682
$ws.Range("K53").Value = 12682
$ws.Range("M53").Value = -12164
$ws.Range("H93").Value = 1467.6
$ws.Range("I93").Value = 1198.1428
$ws.Range("K93").Value = 1198.1428
$ws.Range("M93").Value = 49.85719999999992
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 79857.5
$ws.Range("J125").Value = 79857.5
$ws.Range("L125").Value = 79857.5
$ws.Range("N125").Value = -89697.5
$ws.Range("H132").Value = 950
$ws.Range("I132").Value = 950
$ws.Range("K132").Value = 2850
$ws.Range("M132").Value = -320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 29500
$ws.Range("I50").Value = 29500
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 29500
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -28869
$ws.Range("N50").ClearContents()
$ws.Range("H132").Value = 2846
$ws.Range("I132").Value = 2205.9285
$ws.Range("K132").Value = 6617.7855
$ws.Range("M132").Value = -4087.7855
$ws.Range("H136").Value = 1890.3
$ws.Range("I136").Value = 986.1429000000001
$ws.Range("K136").Value = 2958.4287
$ws.Range("M136").Value = -408.4287000000004
